$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.149.47"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.179.78"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.19"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.88"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.177.65"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("E13").Value = "  +19.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.15"
$ws.Range("E14").Value = "  +7.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.700.73"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.222.79"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.185.70"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("E18").Value = "  +6.50%  "
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.16"
$ws.Range("E20").Value = "  +7.53%  "
$ws.Range("E21").Value = "  +6.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +7.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.46"
$ws.Range("E23").Value = "  +7.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.86"
$ws.Range("E24").Value = "  +3.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.31"
$ws.Range("E25").Value = "  +3.65%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +12.33%  "
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("E29").Value = "  +7.63%  "
$ws.Range("E30").Value = "  +6.77%  "
$ws.Range("E31").Value = "  +13.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +10.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.71"
$ws.Range("E35").Value = "  +7.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.81"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +10.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.74"
$ws.Range("E38").Value = "  +7.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").Value = "  +9.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0423"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.134.76"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("E42").Value = "  +4.70%  "
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("E44").Value = "  +9.90%  "
$ws.Range("E45").Value = "  +13.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.49"
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0607"
$ws.Range("E47").Value = "  +16.71%  "
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").Value = "  +10.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.84"
$ws.Range("E51").Value = "  +2.80%  "
